$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D7").Value = -7.8169
$ws.Range("A9").Value = -21.8093
$ws.Range("D12").Value = -7.057399999999994
$ws.Range("A18").Value = -22.17370000000001
$ws.Range("A20").Value = -19.40079999999999
$ws.Range("D26").Value = -8.517900000000003
$ws.Range("A27").Value = -21.93289999999999
$ws.Range("D27").Value = -8.824199999999998
$ws.Range("D29").Value = -7.269799999999998
$ws.Range("D37").Value = -7.500999999999997
$ws.Range("D38").Value = -8.272199999999996
$ws.Range("D51").Value = -7.830999999999994
$ws.Range("D55").Value = -8.877099999999997
$ws.Range("A69").Value = -21.65919999999999
$ws.Range("D69").Value = -7.097999999999993
$ws.Range("D70").Value = -7.612099999999995
$ws.Range("A76").Value = -19.72449999999998
$ws.Range("A82").Value = -21.86790000000001
$ws.Range("D83").Value = -8.9481
$ws.Range("D102").Value = -7.672099999999995
